$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 45949
$ws.Range("B149").Value = 48654
$ws.Range("E149").Value = 38.26
$ws.Range("F149").Value = -1
$ws.Range("G149").Value = -32.02
$ws.Range("B150").Value = 63902
$ws.Range("E150").Value = 34.04
$ws.Range("F150").Value = 2
$ws.Range("G150").Value = 64.04000000000001
$ws.Range("B161").Value = 53925
$ws.Range("E161").Value = 79.37
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44
$ws.Range("B162").Value = 64350
$ws.Range("E162").Value = 70.63
$ws.Range("F162").Value = 60
$ws.Range("G162").Value = 3986.4
$ws.Range("B163").Value = 57756
$ws.Range("F163").Value = -100
$ws.Range("G163").Value = -6644
$ws.Range("B183").Value = 57552
$ws.Range("E183").Value = 136.86
$ws.Range("F183").Value = -5
$ws.Range("G183").Value = -603.45
$ws.Range("B184").Value = 64329
$ws.Range("E184").Value = 128.32
$ws.Range("F184").Value = 6
$ws.Range("G184").Value = 724.14
$ws.Range("B279").Value = 48706
$ws.Range("E279").Value = 39.8
$ws.Range("F279").Value = -144
$ws.Range("G279").Value = -4795.2
$ws.Range("B280").Value = 64973
$ws.Range("E280").Value = 35.4
$ws.Range("F280").Value = 145
$ws.Range("G280").Value = 4828.5
$ws.Range("B313").Value = 57854
$ws.Range("F313").Value = 2
$ws.Range("G313").Value = 611.6799999999999
$ws.Range("B314").Value = 62997
$ws.Range("F314").Value = 0
$ws.Range("G314").Value = 0
$ws.Range("B316").Value = 61610
$ws.Range("E316").Value = 122.71
$ws.Range("F316").Value = -58
$ws.Range("G316").Value = -5957.18
$ws.Range("B318").Value = 63565
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6
$ws.Range("B351").Value = 57802
$ws.Range("E351").Value = 162.71
$ws.Range("F351").Value = -79
$ws.Range("G351").Value = -11334.92
$ws.Range("B352").Value = 63531
$ws.Range("E352").Value = 152.53
$ws.Range("F352").Value = 80
$ws.Range("G352").Value = 11478.4
$ws.Range("B379").Value = 63564
$ws.Range("F379").Value = 55
$ws.Range("G379").Value = 7095.55
$ws.Range("B380").Value = 65514
$ws.Range("F380").Value = 1
$ws.Range("G380").Value = 129.01
$ws.Range("B389").Value = 57817
$ws.Range("F389").Value = 3
$ws.Range("G389").Value = 239.43
$ws.Range("B390").Value = 62865
$ws.Range("F390").Value = 73
$ws.Range("G390").Value = 5826.13
$ws.Range("B419").Value = 57856
$ws.Range("F419").Value = 2
$ws.Range("G419").Value = 342.66
$ws.Range("B420").Value = 63007
$ws.Range("F420").Value = 898
$ws.Range("G420").Value = 153854.34
$ws.Range("B421").Value = 57857
$ws.Range("F421").Value = 3
$ws.Range("G421").Value = 453.51
$ws.Range("B422").Value = 63008
$ws.Range("F422").Value = 452
$ws.Range("G422").Value = 68328.84
$ws.Range("B457").Value = 31930
$ws.Range("E457").Value = 26.8
$ws.Range("F457").Value = -62
$ws.Range("G457").Value = -1390.04
$ws.Range("B458").Value = 63681
$ws.Range("E458").Value = 23.84
$ws.Range("F458").Value = 43
$ws.Range("G458").Value = 964.0599999999999
$ws.Range("B536").Value = 58047
$ws.Range("D536").Value = 105.54
$ws.Range("E536").Value = 126.1
$ws.Range("F536").Value = 51
$ws.Range("G536").Value = 5382.54
$ws.Range("B537").Value = 47097
$ws.Range("D537").Value = 112.28
$ws.Range("E537").Value = 134.16
$ws.Range("F537").Value = 15
$ws.Range("G537").Value = 1684.2
$ws.Range("B590").Value = 45706
$ws.Range("E590").Value = 23.58
$ws.Range("F590").Value = -202
$ws.Range("G590").Value = -3985.46
$ws.Range("B591").Value = 64922
$ws.Range("E591").Value = 20.98
$ws.Range("F591").Value = 176
$ws.Range("G591").Value = 3472.48
$ws.Range("B593").Value = 64927
$ws.Range("E593").Value = 17.26
$ws.Range("F593").Value = 286
$ws.Range("G593").Value = 4638.92
$ws.Range("B594").Value = 45718
$ws.Range("E594").Value = 19.38
$ws.Range("F594").Value = -294
$ws.Range("G594").Value = -4768.68
$ws.Range("B599").Value = 64925
$ws.Range("E599").Value = 13.97
$ws.Range("F599").Value = 278
$ws.Range("G599").Value = 3655.7
$ws.Range("B600").Value = 45709
$ws.Range("E600").Value = 15.69
$ws.Range("F600").Value = -300
$ws.Range("G600").Value = -3945
$ws.Range("B889").Value = 65079
$ws.Range("F889").Value = 21
$ws.Range("G889").Value = 858.27
$ws.Range("B890").Value = 65362
$ws.Range("F890").Value = 95
$ws.Range("G890").Value = 3882.65
